$d = $word.ActiveDocument

# --- Edit 1: "Specify which technology you are targeting" gains a new
#     trailing run " (Vuforia or AR Foundation)" ---------------------------
$targetText = "Specify which technology you are targeting"
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text.TrimEnd("`r`a") -eq $targetText) {
        $pRange = $para.Range
        # Range covering just the run text, excluding the paragraph mark,
        # so the paragraph's own identity (paraId/rsid/etc.) is untouched.
        $inner = $d.Range($pRange.Start, $pRange.End - 1)
        $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
               '<w:r><w:t>Specify which technology you are targeting</w:t></w:r>' +
               '<w:r><w:t xml:space="preserve"> (Vuforia or AR Foundation)</w:t></w:r>' +
               '</w:p>'
        [void]$inner.InsertXML($xml)
        break
    }
}

# --- Edit 2: fix the "Inlcude" misspelling (and drop its spell-check
#     markup) by rewriting it as three runs "In" + "cl" + "ude the
#     author/creator" ------------------------------------------------------
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text.TrimEnd("`r`a") -eq "Inlcude the author/creator") {
        $full = $para.Range
        $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" ' +
               'w14:paraId="20AECB14" w14:textId="3F1E0687" w:rsidR="008A1F76" ' +
               'w:rsidRDefault="008A1F76" w:rsidP="008A1F76" ' +
               'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
               '<w:pPr><w:pStyle w:val="ListParagraph"/>' +
               '<w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
               '<w:r><w:t>In</w:t></w:r>' +
               '<w:r><w:t>cl</w:t></w:r>' +
               '<w:r><w:t>ude the author/creator</w:t></w:r>' +
               '</w:p>'
        [void]$full.InsertXML($xml)
        break
    }
}
